$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in C8: "dropdownLanguage" -> "dropdownLangauge"
$ws.Range("C8").Value = "dropdownLangauge"

# Lowercase the identifier in B10: "Xpath" -> "xpath"
$ws.Range("B10").Value = "xpath"

# Update the active selection to B10 (was B12)
$ws.Range("B10").Select()
